# Generate Report for Handoff
# Adds two newly-handed-off files (29a78cd5-... and fcce6a46-...) to the
# localization status report: one new row per file on the "Overview" sheet,
# and one new row per file on each language sheet ("zh-cn", "de-de").
#
# Note: a leading "'" forces literal text so Excel doesn't auto-coerce
# "True"/"False" into booleans or collapse an intentional empty string
# into a truly-blank cell.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 4 - 29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md
$ov.Cells.Item(4,1).Value = "29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md"
$ov.Hyperlinks.Add($ov.Cells.Item(4,2), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29a78cd5-handoff/e2e/29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md", "", "", "e2e\29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md")
$ov.Cells.Item(4,3).Value = ".md"
$ov.Cells.Item(4,4).Value = "'"
$ov.Cells.Item(4,5).Value = "Ready for handoff"
$ov.Cells.Item(4,6).Value = "Ready for handoff"
$ov.Cells.Item(4,7).Value = "2016-08-12 14:47:46"
$ov.Cells.Item(4,7).NumberFormat = $dateFmt

# Row 5 - fcce6a46-a380-404a-8ed9-3cede8531aba.md
$ov.Cells.Item(5,1).Value = "fcce6a46-a380-404a-8ed9-3cede8531aba.md"
$ov.Hyperlinks.Add($ov.Cells.Item(5,2), "https://github.com/OpenLocalizationTestOrg/oltest/blob/fcce6a46-handoff/e2e/fcce6a46-a380-404a-8ed9-3cede8531aba.md", "", "", "e2e\fcce6a46-a380-404a-8ed9-3cede8531aba.md")
$ov.Cells.Item(5,3).Value = ".md"
$ov.Cells.Item(5,4).Value = "'"
$ov.Cells.Item(5,5).Value = "Ready for handoff"
$ov.Cells.Item(5,6).Value = "Ready for handoff"
$ov.Cells.Item(5,7).Value = "2016-08-12 14:47:46"
$ov.Cells.Item(5,7).NumberFormat = $dateFmt

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 4 - 29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md
$zh.Hyperlinks.Add($zh.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29a78cd5-handoff/e2e/29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md", "", "", "29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md")
$zh.Cells.Item(4,2).Value = ".md"
$zh.Cells.Item(4,3).Value = "Ready for handoff"
$zh.Cells.Item(4,4).Value = "e2e"
$zh.Cells.Item(4,5).Value = "ht"
$zh.Cells.Item(4,6).Value = "'False"
$zh.Cells.Item(4,7).Value = "29a78cd5-4fc4-41d1-a69d-f17d0bc15365.25ff05c36f7abd50e1c13a7017854b99e45ede8c.zh-cn.xlf"
$zh.Cells.Item(4,8).Value = "2016-08-12 14:47:37"
$zh.Cells.Item(4,8).NumberFormat = $dateFmt
$zh.Cells.Item(4,9).Value = "'"
$zh.Cells.Item(4,10).Value = "'"
$zh.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(4,11).NumberFormat = $dateFmt
$zh.Cells.Item(4,12).Value = "'"
$zh.Cells.Item(4,13).Value = "'True"
$zh.Cells.Item(4,14).Value = "'"
$zh.Cells.Item(4,15).Value = "'False"
$zh.Cells.Item(4,16).Value = "'"

# Row 5 - fcce6a46-a380-404a-8ed9-3cede8531aba.md
$zh.Hyperlinks.Add($zh.Cells.Item(5,1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/fcce6a46-handoff/e2e/fcce6a46-a380-404a-8ed9-3cede8531aba.md", "", "", "fcce6a46-a380-404a-8ed9-3cede8531aba.md")
$zh.Cells.Item(5,2).Value = ".md"
$zh.Cells.Item(5,3).Value = "Ready for handoff"
$zh.Cells.Item(5,4).Value = "e2e"
$zh.Cells.Item(5,5).Value = "ht"
$zh.Cells.Item(5,6).Value = "'False"
$zh.Cells.Item(5,7).Value = "fcce6a46-a380-404a-8ed9-3cede8531aba.4ff472c14935f3f5aa3baf2ae64b806b50c77ff7.zh-cn.xlf"
$zh.Cells.Item(5,8).Value = "2016-08-12 14:47:37"
$zh.Cells.Item(5,8).NumberFormat = $dateFmt
$zh.Cells.Item(5,9).Value = "'"
$zh.Cells.Item(5,10).Value = "'"
$zh.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(5,11).NumberFormat = $dateFmt
$zh.Cells.Item(5,12).Value = "'"
$zh.Cells.Item(5,13).Value = "'True"
$zh.Cells.Item(5,14).Value = "'"
$zh.Cells.Item(5,15).Value = "'False"
$zh.Cells.Item(5,16).Value = "'"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 4 - 29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md
$de.Hyperlinks.Add($de.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29a78cd5-handoff/e2e/29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md", "", "", "29a78cd5-4fc4-41d1-a69d-f17d0bc15365.md")
$de.Cells.Item(4,2).Value = ".md"
$de.Cells.Item(4,3).Value = "Ready for handoff"
$de.Cells.Item(4,4).Value = "e2e"
$de.Cells.Item(4,5).Value = "ht"
$de.Cells.Item(4,6).Value = "'False"
$de.Cells.Item(4,7).Value = "29a78cd5-4fc4-41d1-a69d-f17d0bc15365.25ff05c36f7abd50e1c13a7017854b99e45ede8c.de-de.xlf"
$de.Cells.Item(4,8).Value = "2016-08-12 14:47:46"
$de.Cells.Item(4,8).NumberFormat = $dateFmt
$de.Cells.Item(4,9).Value = "'"
$de.Cells.Item(4,10).Value = "'"
$de.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$de.Cells.Item(4,11).NumberFormat = $dateFmt
$de.Cells.Item(4,12).Value = "'"
$de.Cells.Item(4,13).Value = "'True"
$de.Cells.Item(4,14).Value = "'"
$de.Cells.Item(4,15).Value = "'False"
$de.Cells.Item(4,16).Value = "'"

# Row 5 - fcce6a46-a380-404a-8ed9-3cede8531aba.md
$de.Hyperlinks.Add($de.Cells.Item(5,1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/fcce6a46-handoff/e2e/fcce6a46-a380-404a-8ed9-3cede8531aba.md", "", "", "fcce6a46-a380-404a-8ed9-3cede8531aba.md")
$de.Cells.Item(5,2).Value = ".md"
$de.Cells.Item(5,3).Value = "Ready for handoff"
$de.Cells.Item(5,4).Value = "e2e"
$de.Cells.Item(5,5).Value = "ht"
$de.Cells.Item(5,6).Value = "'False"
$de.Cells.Item(5,7).Value = "fcce6a46-a380-404a-8ed9-3cede8531aba.4ff472c14935f3f5aa3baf2ae64b806b50c77ff7.de-de.xlf"
$de.Cells.Item(5,8).Value = "2016-08-12 14:47:46"
$de.Cells.Item(5,8).NumberFormat = $dateFmt
$de.Cells.Item(5,9).Value = "'"
$de.Cells.Item(5,10).Value = "'"
$de.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$de.Cells.Item(5,11).NumberFormat = $dateFmt
$de.Cells.Item(5,12).Value = "'"
$de.Cells.Item(5,13).Value = "'True"
$de.Cells.Item(5,14).Value = "'"
$de.Cells.Item(5,15).Value = "'False"
$de.Cells.Item(5,16).Value = "'"
